# Apply the authored changes to TestData.xlsx:
#  - Update pCloudy endpoint URL on the "Capabilities" sheet
#  - Swap the sample device / OS data on the "DeviceList" sheet from an
#    iOS device pair to an Android device pair
#  - Move the active/selected sheet + cell-selection state from
#    "Capabilities" (D7) to "DeviceList" (now tabSelected, cell C12)

$wb = $excel.ActiveWorkbook

$capabilities = $wb.Worksheets.Item("Capabilities")
$deviceList   = $wb.Worksheets.Item("DeviceList")

# --- Capabilities sheet: EndPoint values -------------------------------
$capabilities.Range("D2").Value = "https://device.pcloudy.com"
$capabilities.Range("D3").Value = "https://device.pcloudy.com"

# --- DeviceList sheet: swap sample device from iOS to Android ----------
$deviceList.Range("B1").Value = "SAMSUNG_GalaxyA31_Android_11.0.0_32c0a"
$deviceList.Range("C1").Value = "SAMSUNG_GalaxyA50_Android_11.0.0_310bf"
$deviceList.Range("B2").Value = "11.0.0"
$deviceList.Range("C2").Value = "11.0.0"
$deviceList.Range("B3").Value = "pCloudyAndroid"
$deviceList.Range("C3").Value = "pCloudyAndroid"

# --- View state: Capabilities' prior selection moves off D7 to D8, -----
# --- then DeviceList becomes the active/selected sheet at C12 ----------
$null = $capabilities.Activate()
$null = $capabilities.Range("D8").Select()

$null = $deviceList.Activate()
$null = $deviceList.Range("C12").Select()
